$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the "K Noord-West" row entirely; everything below shifts up one row ---
$ws.Rows.Item(11).Delete()

# --- Re-order the "K" group (rows 11-13) so "K Zuid" leads ---
$ws.Range("A11").Value = "K Zuid"
$ws.Range("A12").Value = "K Buitenveldert - Zuidas"
$ws.Range("A13").Value = "K De Pijp - Rivierenbuurt"

# --- Selection left on the active cell / whole row as Excel leaves it after the edits ---
$ws.Range("A14:XFD14").Select()
